$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter period headers (row 8 and row 24): drop the oldest
# quarter ("Q3 ending 1399/06") and append the new quarter
# ("Q1 ending 1401/12") at the end of the period range (columns E:N). ---
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل اول منتهی به 1401/12"

$ws.Range("E24").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F24").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G24").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H24").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I24").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J24").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K24").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L24").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M24").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N24").Value = "فصل اول منتهی به 1401/12"

# --- Shift each metric's quarterly figures left by one quarter and
# append the newly reported figure for the new quarter in column N. ---
$vals = @(41756, 25198, 23647, 43181, 41214, 38152, 36562, 53993, 51203, 38769)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(10, 5 + $i).Value = $vals[$i]
}

$vals = @(-39123, 122582, 86103, 45981, 183296, 209315, 139606, 229624, -58408, 156558)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(13, 5 + $i).Value = $vals[$i]
}

$vals = @(10794, 15915, 7632, 14923, 17256, 22601, 20881, 16508, 20199, 25521)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(14, 5 + $i).Value = $vals[$i]
}

$vals = @(-4434, 1089, 2506, 3302, -6897, 0, 0, 16224, -16224, 4845)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(15, 5 + $i).Value = $vals[$i]
}

$vals = @(9809, 12205, 15355, 23120, 19463, 23281, 23873, 26067, 31762, 44136)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(16, 5 + $i).Value = $vals[$i]
}

$vals = @(383415, 423318, 598338, 456562, 664617, 653400, 1053117, 1036236, 778270, 1210775)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(17, 5 + $i).Value = $vals[$i]
}

$vals = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(18, 5 + $i).Value = $vals[$i]
}

$vals = @(389629, 149007, 186582, 209400, 287646, 210592, 184222, 131291, 1117867, 292068)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(19, 5 + $i).Value = $vals[$i]
}

$vals = @(791846, 749314, 920163, 796469, 1206595, 1157341, 1458261, 1509943, 1924669, 1772672)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(20, 5 + $i).Value = $vals[$i]
}

$vals = @(1039, 1062, 1040, 1093, 1016, 1095, 1106, 1166, 1213, 1256)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(26, 5 + $i).Value = $vals[$i]
}

$vals = @(276, 276, 290, 283, 322, 332, 342, 342, 345, 352)
for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item(27, 5 + $i).Value = $vals[$i]
}

